# Error Calculations and Plots
# Remove the "RM 232" row and the "SC 92" row from the data table;
# remaining rows shift up to close the gaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 currently holds "RM 232" - delete it entirely; rows below shift up.
$ws.Rows(26).Delete()

# After the shift above, the former "SC 92" row (was row 28) is now row 27.
# Delete it too so everything below shifts up again.
$ws.Rows(27).Delete()
